# Update cryptocurrency "Price" values in column D on Sheet1.
# Values are stored as plain text (inline strings) in the workbook, so we
# force the cells to text format before assigning the new value. This keeps
# exact string formatting (e.g. trailing zeros such as "0.01360") intact
# instead of letting Excel coerce the value into a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    2  = "268.42"
    3  = "21.54"
    5  = "0.06207"
    6  = "3.572"
    7  = "6.534"
    8  = "1.396"
    9  = "0.8242"
    10 = "0.1637"
    11 = "0.08192"
    12 = "0.03557"
    13 = "0.03183"
    14 = "0.09198"
    15 = "3.776"
    16 = "0.001632"
    17 = "0.04678"
    18 = "0.006431"
    23 = "2.248"
    24 = "0.01360"
    26 = "0.1243"
    40 = "0.04713"
    42 = "0.003542"
    43 = "0.1122"
    44 = "0.01176"
    45 = "0.00006296"
    46 = "0.0009907"
    47 = "0.00000000750"
    49 = "0.002347"
}

foreach ($row in $updates.Keys) {
    $cell = $ws.Range("D$row")
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$row]
}
